# Update "Tasa Prime - dolar 2021 - Diaria" with new daily rate rows.
# Adds 10 new dated rows (21-09-2021 .. 04-10-2021) to columns A/B, all with
# rate 3.25, continuing directly after the existing last row (181).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dates = @(
    "21-09-2021",
    "22-09-2021",
    "23-09-2021",
    "24-09-2021",
    "27-09-2021",
    "28-09-2021",
    "29-09-2021",
    "30-09-2021",
    "01-10-2021",
    "04-10-2021"
)
$rate = 3.25

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }
$startRow = $lastRow + 1

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $dateCell = $ws.Cells.Item($row, 1)

    # Write the date as a formula that evaluates to the literal text, then
    # convert it to a plain value in place. This guarantees the date string
    # is stored as text (matching the source data) instead of being
    # auto-recognized and coerced into a numeric date serial, while also
    # avoiding the creation of any new/extra cell style.
    $dateCell.Formula = '="' + $dates[$i] + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = $rate
}

$excel.CutCopyMode = 0
